$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 was a category separator row with only column A filled ("NA").
# Fill in the remaining columns with "NA" as well.
$ws.Cells.Item(10, 2).Value = "NA"
$ws.Cells.Item(10, 3).Value = "NA"
$ws.Cells.Item(10, 4).Value = "NA"

# 100 new magic-item rows appended below the existing data (rows 101-200),
# including one more "NA" separator row (row 137) in the middle.
$newItems = @(
    @("Shirt",1,"Utility","Common"),
    @("Shirt",2,"Combat","Uncommon"),
    @("Shirt",1,"Utility","Common"),
    @("Shirt",1,"Utility","Common"),
    @("Shirt",1,"Utility","Uncommon"),
    @("Shirt",3,"Utility","Very Rare"),
    @("Shirt",1,"Utility","Rare"),
    @("Shirt",1,"Utility","Uncommon"),
    @("Shirt",1,"Utility","Very Rare"),
    @("Shirt",1,"Utility","Rare"),
    @("Pants",1,"Utility","Rare"),
    @("Pants",1,"Utility","Rare"),
    @("Pants",1,"Utility","Uncommon"),
    @("Pants",1,"Utility","Rare"),
    @("Pants",2,"Utility","Rare"),
    @("Pants",1,"Utility","Uncommon"),
    @("Pants",1,"Combat","Uncommon"),
    @("Pants",1,"Utility","Uncommon"),
    @("Pants",2,"Combat","Uncommon"),
    @("Pants",1,"Utility","Common"),
    @("Belt",1,"Utility","Very Rare"),
    @("Belt",1,"Utility","Uncommon"),
    @("Belt",1,"Utility","Uncommon"),
    @("Belt",1,"Utility","Uncommon"),
    @("Belt",1,"Utility","Very Rare"),
    @("Belt",2,"Cursed","Rare"),
    @("Belt",1,"Utility","Uncommon"),
    @("Belt",1,"Combat","Very Rare"),
    @("Belt",1,"Utility","Uncommon"),
    @("Belt",1,"Utility","Uncommon"),
    @("Shoes",1,"Utility","Rare"),
    @("Shoes",1,"Utility","Uncommon"),
    @("Shoes",1,"Utility","Uncommon"),
    @("Shoes",2,"Utility","Uncommon"),
    @("Shoes",1,"Utility","Rare"),
    @("Shoes",1,"Utility","Common"),
    @("NA","NA","NA","NA"),
    @("Shoes",3,"Utility","Uncommon"),
    @("Shoes",1,"Utility","Rare"),
    @("Shoes",1,"Cursed","Uncommon"),
    @("Shoes",1,"Utility","Uncommon"),
    @("Shoes",1,"Utility","Uncommon"),
    @("Shoes",1,"Utility","Uncommon"),
    @("Shoes",1,"Utility","Very Rare"),
    @("Shoes",1,"Utility","Rare"),
    @("Shoes",1,"Utility","Rare"),
    @("Shoes",2,"Utility","Uncommon"),
    @("Shoes",2,"Utility","Rare"),
    @("Shoes",1,"Utility","Rare"),
    @("Shoes",1,"Utility","Rare"),
    @("Cloak",1,"Utility","Rare"),
    @("Cloak",3,"Utility","Uncommon"),
    @("Cloak",2,"Utility","Rare"),
    @("Cloak",2,"Combat","Uncommon"),
    @("Cloak",1,"Utility","Very Rare"),
    @("Cloak",2,"Combat","Rare"),
    @("Cloak",1,"Utility","Uncommon"),
    @("Cloak",1,"Utility","Very Rare"),
    @("Cloak",1,"Utility","Rare"),
    @("Cloak",1,"Utility","Uncommon"),
    @("Cloak",1,"Utility","Uncommon"),
    @("Cloak",2,"Utility","Uncommon"),
    @("Cloak",1,"Combat","Uncommon"),
    @("Cloak",1,"Utility","Very Rare"),
    @("Cloak",1,"Combat","Legendary"),
    @("Cloak",1,"Utility","Rare"),
    @("Cloak",3,"Combat","Uncommon"),
    @("Cloak",1,"Utility","Rare"),
    @("Cloak",1,"Combat","Rare"),
    @("Cloak",1,"Combat","Uncommon"),
    @("Hat",1,"Combat","Rare"),
    @("Hat",1,"Utility","Common"),
    @("Hat",1,"Utility","Common"),
    @("Hat",1,"Utility","Uncommon"),
    @("Hat",1,"Utility","Rare"),
    @("Hat",1,"Utility","Uncommon"),
    @("Hat",2,"Utility","Uncommon"),
    @("Hat",1,"Utility","Common"),
    @("Hat",1,"Utility","Uncommon"),
    @("Hat",2,"Utility","Uncommon"),
    @("Gloves",1,"Utility","Uncommon"),
    @("Gloves",1,"Utility","Uncommon"),
    @("Gloves",1,"Combat","Rare"),
    @("Gloves",1,"Combat","Rare"),
    @("Gloves",2,"Combat","Rare"),
    @("Gloves",1,"Utility","Uncommon"),
    @("Gloves",1,"Utility","Rare"),
    @("Gloves",1,"Utility","Uncommon"),
    @("Gloves",1,"Utility","Uncommon"),
    @("Gloves",2,"Combat","Very Rare"),
    @("Scarf",1,"Combat","Rare"),
    @("Cowl",1,"Utility","Uncommon"),
    @("Chaps",1,"Utility","Uncommon"),
    @("Apron",1,"Utility","Uncommon"),
    @("Goggles",1,"Combat","Rare"),
    @("Overalls",1,"Utility","Uncommon"),
    @("Mask",1,"Utility","Uncommon"),
    @("Mask",1,"Utility","Rare"),
    @("Mask",1,"Utility","Legendary"),
    @("Mask",1,"Utility","Uncommon")
)

$startRow = 101
for ($i = 0; $i -lt $newItems.Length; $i++) {
    $r = $startRow + $i
    $item = $newItems[$i]
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
}

# Update the view to match where the author ended up after adding the rows.
$null = $ws.Range("H186").Select()

# Match the page orientation explicitly set on save.
$ws.PageSetup.Orientation = 1 | Out-Null
